$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: maxHP -> stamina (column D)
$ws.Range("D1").Value = "stamina"

# Update stamina values in column D
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 99
$ws.Range("D4").Value = 255

# Update the active selection to D4
$ws.Range("D4").Select()
